$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 2 (H) values updated for logged Week 15 / simulated Week 16
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 217
$wsOff.Range("C2").Value = 150
$wsOff.Range("D2").Value = 52
$wsOff.Range("E2").Value = 17
$wsOff.Range("F2").Value = 3

# Sheet "DEF" - row 2 (H) values updated for logged Week 15 / simulated Week 16
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 196
$wsDef.Range("C2").Value = 119
$wsDef.Range("D2").Value = 31
$wsDef.Range("E2").Value = 9
$wsDef.Range("F2").Value = 3
